# Update the date and the 25 multiplication problems in the document.
$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-01 Sunday", "2024-12-02 Monday"),
    @("348×6=", "157×5="),
    @("172×7=", "630×5="),
    @("822×4=", "333×7="),
    @("210×3=", "326×3="),
    @("277×3=", "392×6="),
    @("224×7=", "375×8="),
    @("500×7=", "334×6="),
    @("869×4=", "174×8="),
    @("544×5=", "739×7="),
    @("415×3=", "932×9="),
    @("263×2=", "988×6="),
    @("534×7=", "117×6="),
    @("890×4=", "983×8="),
    @("300×5=", "914×9="),
    @("309×7=", "544×3="),
    @("548×7=", "204×6="),
    @("798×9=", "495×3="),
    @("833×6=", "102×2="),
    @("698×9=", "901×9="),
    @("394×8=", "815×3="),
    @("651×7=", "997×6="),
    @("841×4=", "238×3="),
    @("324×6=", "214×7="),
    @("389×2=", "874×9="),
    @("518×6=", "678×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Host "Replacements applied: $($replacements.Count)"
